$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Hydrogen
# B3 gets a corrected (larger) value
$ws.Range("B3").Value = 14685408.22818369
# D3 is cleared out entirely (no longer has a numeric value)
$ws.Range("D3").ClearContents()

# Row 7 is relabeled from "Other" to "Biogas" and its D value is corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 8729.165654598692

# A new row 8 is added, re-using the "Other" label that used to live on row 7,
# with a freshly computed D value. Copy A7's formatting (bold/border/centered
# style) down onto the new label cell first, then fill in the values.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 7036.312299965055
